# Updated symbol list with refreshed crypto price/volume/hour data.
# Columns D (Price), E (Volume 1h %) and G (Hora) are stored as literal
# text in the sheet, so NumberFormat is forced to "@" (Text) before each
# write to keep numeric-looking strings (e.g. "305.75", "1.76%", "4")
# from being auto-coerced into numbers/percentages by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "305.75" },
    @{ Cell = "E2"; Value = "1.76%" },
    @{ Cell = "G2"; Value = "4" },
    @{ Cell = "D3"; Value = "36.20" },
    @{ Cell = "E3"; Value = "3.19%" },
    @{ Cell = "G3"; Value = "4" },
    @{ Cell = "D4"; Value = "5.109" },
    @{ Cell = "E4"; Value = "2.36%" },
    @{ Cell = "G4"; Value = "4" },
    @{ Cell = "D5"; Value = "0.08114" },
    @{ Cell = "E5"; Value = "2.70%" },
    @{ Cell = "G5"; Value = "4" },
    @{ Cell = "D6"; Value = "1.942" },
    @{ Cell = "E6"; Value = "0.92%" },
    @{ Cell = "G6"; Value = "4" },
    @{ Cell = "D7"; Value = "7.769" },
    @{ Cell = "E7"; Value = "0.38%" },
    @{ Cell = "G7"; Value = "4" },
    @{ Cell = "D8"; Value = "0.9337" },
    @{ Cell = "E8"; Value = "1.20%" },
    @{ Cell = "G8"; Value = "4" },
    @{ Cell = "D9"; Value = "0.1398" },
    @{ Cell = "E9"; Value = "21.43%" },
    @{ Cell = "G9"; Value = "4" },
    @{ Cell = "D10"; Value = "0.1921" },
    @{ Cell = "E10"; Value = "5.25%" },
    @{ Cell = "G10"; Value = "4" },
    @{ Cell = "D11"; Value = "0.09214" },
    @{ Cell = "E11"; Value = "0.29%" },
    @{ Cell = "G11"; Value = "4" },
    @{ Cell = "D12"; Value = "0.03542" },
    @{ Cell = "E12"; Value = "0.39%" },
    @{ Cell = "G12"; Value = "4" },
    @{ Cell = "D13"; Value = "0.09846" },
    @{ Cell = "E13"; Value = "-0.44%" },
    @{ Cell = "G13"; Value = "4" },
    @{ Cell = "D14"; Value = "0.001418" },
    @{ Cell = "E14"; Value = "1.30%" },
    @{ Cell = "G14"; Value = "4" },
    @{ Cell = "D15"; Value = "0.005817" },
    @{ Cell = "E15"; Value = "0.02%" },
    @{ Cell = "G15"; Value = "4" },
    @{ Cell = "D16"; Value = "3.597" },
    @{ Cell = "E16"; Value = "2.56%" },
    @{ Cell = "G16"; Value = "4" },
    @{ Cell = "D17"; Value = "4.200" },
    @{ Cell = "E17"; Value = "4.74%" },
    @{ Cell = "G17"; Value = "4" },
    @{ Cell = "E18"; Value = "1.87%" },
    @{ Cell = "G18"; Value = "4" },
    @{ Cell = "E19"; Value = "-0.05%" },
    @{ Cell = "G19"; Value = "4" },
    @{ Cell = "E20"; Value = "3.22%" },
    @{ Cell = "G20"; Value = "4" },
    @{ Cell = "D21"; Value = "4.895" },
    @{ Cell = "E21"; Value = "-3.61%" },
    @{ Cell = "G21"; Value = "4" },
    @{ Cell = "D22"; Value = "0.2412" },
    @{ Cell = "E22"; Value = "0.58%" },
    @{ Cell = "G22"; Value = "4" },
    @{ Cell = "D23"; Value = "0.04512" },
    @{ Cell = "E23"; Value = "0.27%" },
    @{ Cell = "G23"; Value = "4" },
    @{ Cell = "E24"; Value = "0.06%" },
    @{ Cell = "G24"; Value = "4" },
    @{ Cell = "D25"; Value = "0.004874" },
    @{ Cell = "E25"; Value = "6.53%" },
    @{ Cell = "G25"; Value = "4" },
    @{ Cell = "E26"; Value = "-0.58%" },
    @{ Cell = "G26"; Value = "4" },
    @{ Cell = "G27"; Value = "4" },
    @{ Cell = "G28"; Value = "4" },
    @{ Cell = "G29"; Value = "4" },
    @{ Cell = "G30"; Value = "4" },
    @{ Cell = "G31"; Value = "4" },
    @{ Cell = "G32"; Value = "4" },
    @{ Cell = "G33"; Value = "4" },
    @{ Cell = "G34"; Value = "4" },
    @{ Cell = "G35"; Value = "4" },
    @{ Cell = "G36"; Value = "4" },
    @{ Cell = "G37"; Value = "4" },
    @{ Cell = "G38"; Value = "4" },
    @{ Cell = "D39"; Value = "0.02011" },
    @{ Cell = "E39"; Value = "5.91%" },
    @{ Cell = "G39"; Value = "4" },
    @{ Cell = "D40"; Value = "0.04929" },
    @{ Cell = "E40"; Value = "5.02%" },
    @{ Cell = "G40"; Value = "4" },
    @{ Cell = "D41"; Value = "0.01100" },
    @{ Cell = "E41"; Value = "15.01%" },
    @{ Cell = "G41"; Value = "4" },
    @{ Cell = "D42"; Value = "0.007657" },
    @{ Cell = "E42"; Value = "0.56%" },
    @{ Cell = "G42"; Value = "4" },
    @{ Cell = "D43"; Value = "0.1381" },
    @{ Cell = "E43"; Value = "4.35%" },
    @{ Cell = "G43"; Value = "4" },
    @{ Cell = "D44"; Value = "0.002104" },
    @{ Cell = "E44"; Value = "-0.26%" },
    @{ Cell = "G44"; Value = "4" },
    @{ Cell = "D45"; Value = "0.01058" },
    @{ Cell = "E45"; Value = "-4.75%" },
    @{ Cell = "G45"; Value = "4" },
    @{ Cell = "D46"; Value = "0.00006453" },
    @{ Cell = "E46"; Value = "7.46%" },
    @{ Cell = "G46"; Value = "4" },
    @{ Cell = "E47"; Value = "0.19%" },
    @{ Cell = "G47"; Value = "4" },
    @{ Cell = "G48"; Value = "4" },
    @{ Cell = "E49"; Value = "-8.58%" },
    @{ Cell = "G49"; Value = "4" },
    @{ Cell = "D50"; Value = "0.00002104" },
    @{ Cell = "E50"; Value = "0.19%" },
    @{ Cell = "G50"; Value = "4" },
    @{ Cell = "D51"; Value = "0.0002003" },
    @{ Cell = "E51"; Value = "0.19%" },
    @{ Cell = "G51"; Value = "4" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
